# Updates cryptocurrency price/volume data to match latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.382.67'
$ws.Range('E2').Value = '  +2.97%  '
$ws.Range('D3').Value = '1.915.27'
$ws.Range('E3').Value = '  +1.47%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''248.80'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.00%  '
$ws.Range('D6').Value = '''0.694'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.56%  '
$ws.Range('D7').Value = '''0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '''43.98'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.25%  '
$ws.Range('D9').Value = '''58.32'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +7.87%  '
$ws.Range('D10').Value = '''0.364'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.19%  '
$ws.Range('E11').Value = '  +2.69%  '
$ws.Range('D12').Value = '''0.0992'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.12%  '
$ws.Range('D13').Value = '''14.43'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +8.98%  '
$ws.Range('E14').Value = '  +5.59%  '
$ws.Range('D15').Value = '2.194.73'
$ws.Range('E15').Value = '  +1.57%  '
$ws.Range('D16').Value = '''5.11'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.51%  '
$ws.Range('D17').Value = '1.916.90'
$ws.Range('E17').Value = '  +1.63%  '
$ws.Range('D18').Value = '36.332.20'
$ws.Range('E18').Value = '  +2.64%  '
$ws.Range('D19').Value = '''74.29'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.73%  '
$ws.Range('D20').Value = '0.0₃0847'
$ws.Range('E20').Value = '  +3.03%  '
$ws.Range('D21').Value = '''251.85'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.02%  '
$ws.Range('D22').Value = '''13.14'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.89%  '
$ws.Range('D23').Value = '''5.15'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.27%  '
$ws.Range('E24').Value = '  +0.17%  '
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('D26').Value = '''2.19'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.28%  '
$ws.Range('D27').Value = '''167.74'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.86%  '
$ws.Range('D28').Value = '''8.75'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.94%  '
$ws.Range('D29').Value = '''18.81'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.73%  '
$ws.Range('E30').Value = '  +1.24%  '
$ws.Range('E31').Value = '  +6.15%  '
$ws.Range('E32').Value = '  +4.18%  '
$ws.Range('D33').Value = '''1.96'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +7.06%  '
$ws.Range('D34').Value = '''4.32'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.55%  '
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('D36').Value = '''0.0847'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +22.12%  '
$ws.Range('E37').Value = '  -14.81%  '
$ws.Range('E38').Value = '  +1.19%  '
$ws.Range('D39').Value = '''2.00'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.13%  '
$ws.Range('D40').Value = '''105.13'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +8.91%  '
$ws.Range('E41').Value = '  +3.56%  '
$ws.Range('D42').Value = '''15.64'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +27.83%  '
$ws.Range('D43').Value = '''17.06'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.53%  '
$ws.Range('E44').Value = '  +2.44%  '
$ws.Range('D45').Value = '1.340.75'
$ws.Range('E45').Value = '  +3.15%  '
$ws.Range('E46').Value = '  +2.05%  '
$ws.Range('E47').Value = '  +1.59%  '
$ws.Range('E48').Value = '  +0.97%  '
$ws.Range('D49').Value = '''2.79'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.07%  '
$ws.Range('D50').Value = '''6.42'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.82%  '
$ws.Range('D51').Value = '2.093.13'
$ws.Range('E51').Value = '  +1.20%  '
